$wb = $excel.ActiveWorkbook

# --- Update selection on "Tentaive Dates" (sheet1): C5 -> C11 ---
$ws1 = $wb.Worksheets.Item("Tentaive Dates")
$ws1.Range("C11").Select()

# --- Update selection on "Mock interview Schedules" (sheet2): B9 -> B10 ---
$ws2 = $wb.Worksheets.Item("Mock interview Schedules")
$ws2.Activate()
$ws2.Range("B10").Select()

# --- Add the new sheet "Resume_CV Preparation" after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "Resume_CV Preparation"

# Column widths (approx. 26.57 / 27.71 / 36.86 chars; closest achievable via COM pixel rounding)
$ws5.Columns.Item(1).ColumnWidth = 25.666666666666664
$ws5.Columns.Item(2).ColumnWidth = 26.833333333333336
$ws5.Columns.Item(3).ColumnWidth = 36.0

# Row data
$ws5.Range("A1").Value = "Amit"
$ws5.Range("A2").Value = "Pramod"
$ws5.Range("A3").Value = "Divyasree"
$ws5.Range("A4").Value = "Keerthana"
$ws5.Range("A5").Value = "Bhuvaneshwari"
$ws5.Range("A6").Value = "Lohanathan"
$ws5.Range("A7").Value = "Akhila"
$ws5.Range("B7").Value = "Shared and given 1st review"
$ws5.Range("C7").Value = "Review comments Implemented or not"
$ws5.Range("A8").Value = "Dikshith"
$ws5.Range("A9").Value = "Vincy"
$ws5.Range("A10").Value = "Ranjitha"
$ws5.Range("A11").Value = "Sreedhar"
$ws5.Range("A12").Value = "Elavarsan"
$ws5.Range("A13").Value = "Vincy"

$ws5.Range("C7").Select()

# --- Restore "Tentaive Dates" as the active sheet/tab ---
$ws1.Activate()
